# Adds multi-task (multi-file) loop support around the "source code" block:
#   - before the "Исходный код" heading paragraph: insert a
#       "{%p for task in tasks %}" paragraph
#   - turn the old heading's formatting into a new "Задание N. task.name"
#       paragraph (keeps the old pStyle="836" heading look)
#   - the former heading paragraph keeps only its suppressLineNumbers flag
#   - the {{source}} placeholder becomes {{ task.source }}
#   - after it, insert a "{%p endfor %}" paragraph
#   - the following "Тестирование" heading switches from pStyle 836 to 835

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$d = $word.ActiveDocument

function Find-ParaIndex($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Find-ParaIndex: text not found: " + $searchText)
    }
    $pCount = $d.Paragraphs.Count
    for ($i = 1; $i -le $pCount; $i++) {
        $pp = $d.Paragraphs.Item($i)
        if ($pp.Range.Start -le $rng.Start -and $pp.Range.End -ge $rng.End) {
            return $i
        }
    }
    throw ("Find-ParaIndex: containing paragraph not found for: " + $searchText)
}

# --- locate the "Исходный код" heading paragraph via the {{source}} anchor --
$idxSource = Find-ParaIndex("{{source}}")
$idxCode = $idxSource - 1

# 1) Insert the "{%p for task in tasks %}" paragraph right before the
#    "Исходный код" heading paragraph.
$codePara = $d.Paragraphs.Item($idxCode)
$codePara.Range.InsertParagraphBefore() | Out-Null

$forPara = $d.Paragraphs.Item($idxCode)
$forXml = '<w:p xmlns:w="' + $w + '">' +
    '<w:pPr><w:rPr><w:highlight w:val="none"/></w:rPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">{%p for task in tasks %}</w:t></w:r>' +
    '<w:r/>' +
    '</w:p>'
$forPara.Range.InsertXML($forXml) | Out-Null

# indices shift by one after the insertion above
$idxCode = $idxCode + 1
$idxSource = $idxSource + 1

# 2) Insert the "Задание {{loop.index}}. {{ task.name }}" paragraph right
#    before the (shifted) "Исходный код" heading paragraph, reusing the
#    pStyle="836" subsection look the heading used to have.
$codePara = $d.Paragraphs.Item($idxCode)
$codePara.Range.InsertParagraphBefore() | Out-Null

$taskPara = $d.Paragraphs.Item($idxCode)
$taskXml = '<w:p xmlns:w="' + $w + '">' +
    '<w:pPr><w:pStyle w:val="836"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Задание {{loop.index}}</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">. {{ </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">task.name</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> }}</w:t></w:r>' +
    '<w:r/>' +
    '</w:p>'
$taskPara.Range.InsertXML($taskXml) | Out-Null

# indices shift by one again
$idxCode = $idxCode + 1
$idxSource = $idxSource + 1

# 3) Strip the old heading's pPr down to just suppressLineNumbers (its runs
#    / text stay exactly as they were: "Исхо" + "дный код").
$codePara = $d.Paragraphs.Item($idxCode)
$codeXml = '<w:p xmlns:w="' + $w + '">' +
    '<w:pPr><w:suppressLineNumbers w:val="0"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Исхо</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">дный код</w:t></w:r>' +
    '<w:r/>' +
    '</w:p>'
$codePara.Range.InsertXML($codeXml) | Out-Null

# 4) Replace {{source}} with {{ task.source }}, keeping the paragraph's
#    own pPr (pStyle="841", jc="left") untouched.
$sourcePara = $d.Paragraphs.Item($idxSource)
$sourceXml = '<w:p xmlns:w="' + $w + '">' +
    '<w:pPr><w:pStyle w:val="841"/><w:jc w:val="left"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">{{ task.source }}</w:t></w:r>' +
    '<w:r/>' +
    '</w:p>'
$sourcePara.Range.InsertXML($sourceXml) | Out-Null

# 5) Insert the "{%p endfor %}" paragraph right after the {{ task.source }}
#    paragraph.
$sourcePara = $d.Paragraphs.Item($idxSource)
$sourcePara.Range.InsertParagraphAfter() | Out-Null

$endforPara = $d.Paragraphs.Item($idxSource + 1)
$endforXml = '<w:p xmlns:w="' + $w + '">' +
    '<w:r><w:t xml:space="preserve">{%p endfor %}</w:t></w:r>' +
    '<w:r/>' +
    '</w:p>'
$endforPara.Range.InsertXML($endforXml) | Out-Null

# 6) Switch the following "Тестирование" heading from pStyle 836 to 835.
$idxTest = Find-ParaIndex("Тестирование")
$testPara = $d.Paragraphs.Item($idxTest)
$testXml = '<w:p xmlns:w="' + $w + '">' +
    '<w:pPr><w:pStyle w:val="835"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Тестирование</w:t></w:r>' +
    '<w:r/>' +
    '</w:p>'
$testPara.Range.InsertXML($testXml) | Out-Null

Write-Output "done"
